# Auto-generated edit script: updates FFXIV leve market-price data cells
# per the scheduled-runner data refresh (columns H-N on specific rows).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2505465
$ws.Range("J17").Value = 2708522.8
$ws.Range("L17").Value = 8125568.399999999
$ws.Range("N17").Value = -8125904.399999999

$ws.Range("H76").Value = 3135.7058
$ws.Range("I76").Value = 3023.6155
$ws.Range("K76").Value = 3023.6155
$ws.Range("M76").Value = -2708.6155

$ws.Range("H79").Value = 3135.7058
$ws.Range("I79").Value = 3023.6155
$ws.Range("K79").Value = 3023.6155
$ws.Range("M79").Value = -1931.6155

$ws.Range("H113").Value = 37041040
$ws.Range("I113").Value = 52635076
$ws.Range("J113").Value = 5201.5
$ws.Range("K113").Value = 52635076
$ws.Range("L113").Value = 5201.5
$ws.Range("M113").Value = -52631822
$ws.Range("N113").Value = -11709.5

$ws.Range("H131").Value = 1828.7858
$ws.Range("I131").Value = 1296.36
$ws.Range("J131").Value = 2611.7646
$ws.Range("K131").Value = 3889.08
$ws.Range("L131").Value = 7835.293799999999
$ws.Range("M131").Value = 1150.92
$ws.Range("N131").Value = -17915.2938

$ws.Range("H135").Value = 9806727
$ws.Range("I135").Value = 490.33334
$ws.Range("J135").Value = 55569170
$ws.Range("K135").Value = 4413.00006
$ws.Range("L135").Value = 500122530
$ws.Range("M135").Value = -1878.00006
$ws.Range("N135").Value = -500127600

$ws.Range("H137").Value = 1801.8857
$ws.Range("I137").Value = 1809.3928
$ws.Range("K137").Value = 5428.178400000001
$ws.Range("M137").Value = -2878.178400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 300.5
$ws.Range("I4").Value = 301
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 301
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -185
$ws.Range("N4").Value = -532

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H32").Value = 3107.6292
$ws.Range("I32").Value = 2727.9875
$ws.Range("K32").Value = 2727.9875
$ws.Range("M32").Value = -2440.9875

$ws.Range("H37").Value = 15512
$ws.Range("I37").Value = 1034
$ws.Range("J37").Value = 29990
$ws.Range("K37").Value = 1034
$ws.Range("L37").Value = 29990
$ws.Range("M37").Value = -761
$ws.Range("N37").Value = -30536

$ws.Range("H59").Value = 20750
$ws.Range("J59").Value = 20750
$ws.Range("L59").Value = 20750
$ws.Range("N59").Value = -22358

$ws.Range("H60").Value = 22000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 22000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 22000
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -23466

$ws.Range("H63").Value = 2070
$ws.Range("I63").Value = 1748.3334
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 1748.3334
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -1062.3334
$ws.Range("N63").Value = -5372

$ws.Range("H66").Value = 2070
$ws.Range("I66").Value = 1748.3334
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 8741.666999999999
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -5309.666999999999
$ws.Range("N66").Value = -26864

$ws.Range("H74").Value = 58826936
$ws.Range("I74").Value = 71432210
$ws.Range("K74").Value = 71432210
$ws.Range("M74").Value = -71431336

$ws.Range("H77").Value = 58826936
$ws.Range("I77").Value = 71432210
$ws.Range("K77").Value = 357161050
$ws.Range("M77").Value = -357156682

$ws.Range("H122").Value = 1562.0541
$ws.Range("I122").Value = 1151.9
$ws.Range("J122").Value = 3319.8572
$ws.Range("K122").Value = 3455.7
$ws.Range("L122").Value = 9959.571599999999
$ws.Range("M122").Value = -1005.7
$ws.Range("N122").Value = -14859.5716

$ws.Range("H132").Value = 16751.559
$ws.Range("I132").Value = 1967.5
$ws.Range("K132").Value = 5902.5
$ws.Range("M132").Value = -3372.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H82").Value = 8328.4
$ws.Range("I82").Value = 8328.4
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 8328.4
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -7945.4
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 8328.4
$ws.Range("I85").Value = 8328.4
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 8328.4
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -7002.4
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3762.8635
$ws.Range("J31").Value = 5046.9
$ws.Range("L31").Value = 5046.9
$ws.Range("N31").Value = -5636.9

$ws.Range("H34").Value = 3762.8635
$ws.Range("J34").Value = 5046.9
$ws.Range("L34").Value = 5046.9
$ws.Range("N34").Value = -5450.9

$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290

$ws.Range("H132").Value = 2172.5
$ws.Range("I132").Value = 1493.1538
$ws.Range("K132").Value = 4479.4614
$ws.Range("M132").Value = -1949.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 50.65
$ws.Range("I12").Value = 5.714286
$ws.Range("J12").Value = 74.84614999999999
$ws.Range("K12").Value = 17.142858
$ws.Range("L12").Value = 224.53845
$ws.Range("M12").Value = 155.857142
$ws.Range("N12").Value = -570.53845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3442.889
$ws.Range("I80").Value = 2784.8462
$ws.Range("J80").Value = 3814.8262
$ws.Range("K80").Value = 2784.8462
$ws.Range("L80").Value = 3814.8262
$ws.Range("M80").Value = -1786.8462
$ws.Range("N80").Value = -5810.8262

$ws.Range("H83").Value = 3442.889
$ws.Range("I83").Value = 2784.8462
$ws.Range("J83").Value = 3814.8262
$ws.Range("K83").Value = 13924.231
$ws.Range("L83").Value = 19074.131
$ws.Range("M83").Value = -8932.231
$ws.Range("N83").Value = -29058.131

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5147.6523
$ws.Range("I7").Value = 3517.3635
$ws.Range("K7").Value = 3517.3635
$ws.Range("M7").Value = -3405.3635

$ws.Range("H22").Value = 2454.3572
$ws.Range("I22").Value = 1867
$ws.Range("J22").Value = 2894.875
$ws.Range("K22").Value = 1867
$ws.Range("L22").Value = 2894.875
$ws.Range("M22").Value = -1572
$ws.Range("N22").Value = -3484.875

$ws.Range("H27").Value = 2454.3572
$ws.Range("I27").Value = 1867
$ws.Range("J27").Value = 2894.875
$ws.Range("K27").Value = 1867
$ws.Range("L27").Value = 2894.875
$ws.Range("M27").Value = -1760
$ws.Range("N27").Value = -3108.875

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H93").Value = 975.6842
$ws.Range("I93").Value = 1032.5834
$ws.Range("J93").Value = 878.1429000000001
$ws.Range("K93").Value = 1032.5834
$ws.Range("L93").Value = 878.1429000000001
$ws.Range("M93").Value = 215.4166
$ws.Range("N93").Value = -3374.1429

$ws.Range("H126").Value = 5147.6523
$ws.Range("I126").Value = 3517.3635
$ws.Range("K126").Value = 10552.0905
$ws.Range("M126").Value = -8082.0905

$ws.Range("H136").Value = 812.6667
$ws.Range("I136").Value = 812.6667
$ws.Range("K136").Value = 2438.0001
$ws.Range("M136").Value = 111.9998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1060.1111
$ws.Range("I132").Value = 709.3714
$ws.Range("J132").Value = 2287.7
$ws.Range("K132").Value = 2128.1142
$ws.Range("L132").Value = 6863.099999999999
$ws.Range("M132").Value = 401.8858
$ws.Range("N132").Value = -11923.1

$ws.Range("H136").Value = 20835822
$ws.Range("I136").Value = 31251196
$ws.Range("J136").Value = 5075.625
$ws.Range("K136").Value = 93753588
$ws.Range("L136").Value = 15226.875
$ws.Range("M136").Value = -93751038
$ws.Range("N136").Value = -20326.875

